$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 429.75
$ws.Range("I18").Value = 476.66666
$ws.Range("K18").Value = 476.66666
$ws.Range("M18").Value = -192.66666
$ws.Range("H20").Value = 200
$ws.Range("I20").Value = 200
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = 30
$ws.Range("H28").Value = 2333.0715
$ws.Range("I28").Value = 1336.5
$ws.Range("J28").Value = 4824.5
$ws.Range("K28").Value = 1336.5
$ws.Range("L28").Value = 4824.5
$ws.Range("M28").Value = -851.5
$ws.Range("N28").Value = -5794.5
$ws.Range("H35").Value = 200
$ws.Range("I35").Value = 200
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 200
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").Value = 179
$ws.Range("H38").Value = 294
$ws.Range("I38").Value = 294
$ws.Range("K38").Value = 882
$ws.Range("M38").Value = -510
$ws.Range("H51").Value = 12554.667
$ws.Range("I51").Value = 10598.6
$ws.Range("J51").Value = 14999.75
$ws.Range("K51").Value = 10598.6
$ws.Range("L51").Value = 14999.75
$ws.Range("M51").Value = -10114.6
$ws.Range("N51").Value = -15967.75
$ws.Range("H62").Value = 1900
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H64").Value = 7496
$ws.Range("I64").Value = 3992
$ws.Range("J64").Value = 11000
$ws.Range("K64").Value = 3992
$ws.Range("L64").Value = 11000
$ws.Range("M64").Value = -3744
$ws.Range("N64").Value = -11496
$ws.Range("H65").Value = 1900
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H67").Value = 7496
$ws.Range("I67").Value = 3992
$ws.Range("J67").Value = 11000
$ws.Range("K67").Value = 3992
$ws.Range("L67").Value = 11000
$ws.Range("M67").Value = -3134
$ws.Range("N67").Value = -12716
$ws.Range("H98").Value = 2622.8462
$ws.Range("I98").Value = 2220.5833
$ws.Range("K98").Value = 2220.5833
$ws.Range("M98").Value = -722.5832999999998
$ws.Range("H113").Value = 3030.75
$ws.Range("I113").Value = 2695.3125
$ws.Range("K113").Value = 2695.3125
$ws.Range("M113").Value = 558.6875
$ws.Range("H122").Value = 2622.8462
$ws.Range("I122").Value = 2220.5833
$ws.Range("K122").Value = 6661.749899999999
$ws.Range("M122").Value = -4211.749899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 255.2
$ws.Range("I5").Value = 314
$ws.Range("J5").Value = 20
$ws.Range("K5").Value = 314
$ws.Range("L5").Value = 20
$ws.Range("M5").Value = -202
$ws.Range("N5").Value = -244
$ws.Range("H28").Value = 12462.8
$ws.Range("I28").Value = 12462.8
$ws.Range("K28").Value = 12462.8
$ws.Range("M28").Value = -12270.8
$ws.Range("H32").Value = 3088.5925
$ws.Range("I32").Value = 2438.1538
$ws.Range("K32").Value = 2438.1538
$ws.Range("M32").Value = -2151.1538
$ws.Range("H41").Value = 35928
$ws.Range("I41").Value = 9856
$ws.Range("K41").Value = 9856
$ws.Range("M41").Value = -9442
$ws.Range("H88").Value = 2923
$ws.Range("I88").Value = 2909.5
$ws.Range("K88").Value = 2909.5
$ws.Range("M88").Value = -2503.5
$ws.Range("H91").Value = 2923
$ws.Range("I91").Value = 2909.5
$ws.Range("K91").Value = 2909.5
$ws.Range("M91").Value = -1505.5
$ws.Range("H99").Value = 12462.8
$ws.Range("I99").Value = 12462.8
$ws.Range("K99").Value = 12462.8
$ws.Range("M99").Value = -9467.799999999999
$ws.Range("H102").Value = 3020.25
$ws.Range("I102").Value = 3020.25
$ws.Range("K102").Value = 3020.25
$ws.Range("M102").Value = -1398.25
$ws.Range("H122").Value = 6940.2
$ws.Range("I122").Value = 6940.2
$ws.Range("K122").Value = 20820.6
$ws.Range("M122").Value = -18370.6
$ws.Range("H132").Value = 5260.355
$ws.Range("I132").Value = 4321.4546
$ws.Range("J132").Value = 7555.4443
$ws.Range("K132").Value = 12964.3638
$ws.Range("L132").Value = 22666.3329
$ws.Range("M132").Value = -10434.3638
$ws.Range("N132").Value = -27726.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 255.2
$ws.Range("I4").Value = 314
$ws.Range("J4").Value = 20
$ws.Range("K4").Value = 314
$ws.Range("L4").Value = 20
$ws.Range("M4").Value = -199
$ws.Range("N4").Value = -250

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 2499
$ws.Range("I17").Value = 2499
$ws.Range("K17").Value = 2499
$ws.Range("M17").Value = -2325
$ws.Range("H45").Value = 17537
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 25074
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 25074
$ws.Range("M45").Value = -9407
$ws.Range("N45").Value = -26260
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
$ws.Range("H107").Value = 604.7917
$ws.Range("I107").Value = 462.66666
$ws.Range("K107").Value = 462.66666
$ws.Range("M107").Value = 1457.33334
$ws.Range("H122").Value = 1131.8235
$ws.Range("I122").Value = 585.2
$ws.Range("K122").Value = 1755.6
$ws.Range("M122").Value = 694.3999999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 409
$ws.Range("I13").Value = 10.75
$ws.Range("J13").Value = 2002
$ws.Range("K13").Value = 32.25
$ws.Range("L13").Value = 6006
$ws.Range("M13").Value = 135.75
$ws.Range("N13").Value = -6342
$ws.Range("H131").Value = 3962.5
$ws.Range("I131").Value = 3500
$ws.Range("K131").Value = 10500
$ws.Range("M131").Value = -5460
$ws.Range("H132").Value = 1000
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 13899
$ws.Range("I70").Value = 13899
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 13899
$ws.Range("L70").ClearContents()
$ws.Range("M70").Value = -13629
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 13899
$ws.Range("I73").Value = 13899
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 13899
$ws.Range("L73").ClearContents()
$ws.Range("M73").Value = -12963
$ws.Range("N73").Value = 0
$ws.Range("H102").Value = 1208.5714
$ws.Range("I102").Value = 976.6667
$ws.Range("J102").Value = 2600
$ws.Range("K102").Value = 976.6667
$ws.Range("L102").Value = 2600
$ws.Range("M102").Value = 645.3333
$ws.Range("N102").Value = -5844
$ws.Range("H113").Value = 2607.5
$ws.Range("I113").Value = 2560
$ws.Range("K113").Value = 2560
$ws.Range("M113").Value = -390

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 3636
$ws.Range("I4").Value = 3636
$ws.Range("K4").Value = 3636
$ws.Range("M4").Value = -3523
$ws.Range("H22").Value = 2315.1667
$ws.Range("I22").Value = 773.5
$ws.Range("J22").Value = 5398.5
$ws.Range("K22").Value = 773.5
$ws.Range("L22").Value = 5398.5
$ws.Range("M22").Value = -478.5
$ws.Range("N22").Value = -5988.5
$ws.Range("H27").Value = 2315.1667
$ws.Range("I27").Value = 773.5
$ws.Range("J27").Value = 5398.5
$ws.Range("K27").Value = 773.5
$ws.Range("L27").Value = 5398.5
$ws.Range("M27").Value = -666.5
$ws.Range("N27").Value = -5612.5
$ws.Range("H28").Value = 3636
$ws.Range("I28").Value = 3636
$ws.Range("K28").Value = 3636
$ws.Range("M28").Value = -3404
$ws.Range("H37").Value = 3636
$ws.Range("I37").Value = 3636
$ws.Range("K37").Value = 3636
$ws.Range("M37").Value = -3529
$ws.Range("H46").Value = 1179.1852
$ws.Range("J46").Value = 979
$ws.Range("L46").Value = 979
$ws.Range("N46").Value = -1355
$ws.Range("H82").Value = 2224.5
$ws.Range("I82").Value = 2224.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2224.5
$ws.Range("L82").Value = 0
$ws.Range("M82").ClearContents()
$ws.Range("N82").Value = -1863.5
$ws.Range("H85").Value = 2224.5
$ws.Range("I85").Value = 2224.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2224.5
$ws.Range("L85").Value = 0
$ws.Range("M85").ClearContents()
$ws.Range("N85").Value = -976.5
$ws.Range("H93").Value = 3175.889
$ws.Range("I93").Value = 2938.4
$ws.Range("J93").Value = 3472.75
$ws.Range("K93").Value = 2938.4
$ws.Range("L93").Value = 3472.75
$ws.Range("M93").Value = -1690.4
$ws.Range("N93").Value = -5968.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 17000
$ws.Range("I52").Value = 17000
$ws.Range("K52").Value = 17000
$ws.Range("M52").Value = -16774
$ws.Range("H55").Value = 4500
$ws.Range("I55").Value = 1500
$ws.Range("J55").Value = 7500
$ws.Range("K55").Value = 1500
$ws.Range("L55").Value = 7500
$ws.Range("M55").Value = -1223
$ws.Range("N55").Value = -8054
$ws.Range("H122").Value = 2001
$ws.Range("I122").Value = 2001
$ws.Range("K122").Value = 6003
$ws.Range("M122").Value = -3553
$ws.Range("H132").Value = 2162.3635
$ws.Range("I132").Value = 2219.6
$ws.Range("K132").Value = 6658.799999999999
$ws.Range("M132").Value = -4128.799999999999
